$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 467 (existing rows 467.. shift down to 469..)
$ws.Range("A467:A468").EntireRow.Insert()

# New row 467: Crimpson Seedless, Primera
$ws.Cells.Item(467, 1).Value = 5
$ws.Cells.Item(467, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(467, 3).Value = "Maule"
$ws.Cells.Item(467, 4).Value = 44706
$ws.Cells.Item(467, 5).Value = 7
$ws.Cells.Item(467, 6).Value = "Fruta"
$ws.Cells.Item(467, 7).Value = 100109
$ws.Cells.Item(467, 8).Value = "Uva"
$ws.Cells.Item(467, 9).Value = 100109001
$ws.Cells.Item(467, 10).Value = "Uva"
$ws.Cells.Item(467, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(467, 12).Value = "Primera"
$ws.Cells.Item(467, 13).Value = 180
$ws.Cells.Item(467, 14).Value = 12000
$ws.Cells.Item(467, 15).Value = 12000
$ws.Cells.Item(467, 16).Value = 12000
$ws.Cells.Item(467, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(467, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(467, 19).Value = 667
$ws.Cells.Item(467, 20).Value = 18

# New row 468: Red Globe, Primera
$ws.Cells.Item(468, 1).Value = 5
$ws.Cells.Item(468, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(468, 3).Value = "Maule"
$ws.Cells.Item(468, 4).Value = 44706
$ws.Cells.Item(468, 5).Value = 7
$ws.Cells.Item(468, 6).Value = "Fruta"
$ws.Cells.Item(468, 7).Value = 100109
$ws.Cells.Item(468, 8).Value = "Uva"
$ws.Cells.Item(468, 9).Value = 100109001
$ws.Cells.Item(468, 10).Value = "Uva"
$ws.Cells.Item(468, 11).Value = "Red Globe"
$ws.Cells.Item(468, 12).Value = "Primera"
$ws.Cells.Item(468, 13).Value = 200
$ws.Cells.Item(468, 14).Value = 10000
$ws.Cells.Item(468, 15).Value = 10000
$ws.Cells.Item(468, 16).Value = 10000
$ws.Cells.Item(468, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(468, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(468, 19).Value = 556
$ws.Cells.Item(468, 20).Value = 18
